$wb = $excel.ActiveWorkbook

$oldId = "62ac2419-5f30-4111-9d4a-08ed0bb0638d"
$newId = "3daeb6c3-5e79-4bff-a0e7-1d5c3314f394"

# --- Overview sheet ------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newId.md"
$wsOverview.Range("B2").Value = "e2e\$newId.md"
$wsOverview.Range("G2").Value = "2016-08-20 23:01:59"

foreach ($h in $wsOverview.Hyperlinks) {
    $h.TextToDisplay = "e2e\$newId.md"
}

# --- zh-cn sheet -----------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "$newId.md"
$wsZhCn.Range("G2").Value = "$newId.d94331dfbdcbbfb866c52d7af1d3f8745c3d09e4.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-20 23:01:55"

foreach ($h in $wsZhCn.Hyperlinks) {
    $h.TextToDisplay = "$newId.md"
}

# --- de-de sheet -----------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "$newId.md"
$wsDeDe.Range("G2").Value = "$newId.d94331dfbdcbbfb866c52d7af1d3f8745c3d09e4.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-20 23:01:59"

foreach ($h in $wsDeDe.Hyperlinks) {
    $h.TextToDisplay = "$newId.md"
}
